$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BY (col 77) width to match the rest of the data columns (stored width 12)
$ws.Range("BY1").ColumnWidth = 11.166666666666666

# Header cell BY1: date label "2024/11/24" as literal text (not an auto-parsed date),
# matching the style used by the other header cells (BX1 etc.)
$ws.Range("BY1").Value = "'2024/11/24"
$ws.Range("BX1").Copy()
$ws.Range("BY1").PasteSpecial(-4122)

# Data rows 2-53: new values for 2024/11/24, each styled to match the sheet's
# value-based highlighting convention (>=140 normal, 125-139.9 light blue, <125 yellow)
$ws.Range("BY2").Value = 207.1
$ws.Range("B2").Copy()
$ws.Range("BY2").PasteSpecial(-4122)
$ws.Range("BY3").Value = 120.9
$ws.Range("D2").Copy()
$ws.Range("BY3").PasteSpecial(-4122)
$ws.Range("BY4").Value = 165.4
$ws.Range("B2").Copy()
$ws.Range("BY4").PasteSpecial(-4122)
$ws.Range("BY5").Value = 218.6
$ws.Range("B2").Copy()
$ws.Range("BY5").PasteSpecial(-4122)
$ws.Range("BY6").Value = 180
$ws.Range("B2").Copy()
$ws.Range("BY6").PasteSpecial(-4122)
$ws.Range("BY7").Value = 142.4
$ws.Range("B2").Copy()
$ws.Range("BY7").PasteSpecial(-4122)
$ws.Range("BY8").Value = 144.3
$ws.Range("B2").Copy()
$ws.Range("BY8").PasteSpecial(-4122)
$ws.Range("BY9").Value = 130.9
$ws.Range("N2").Copy()
$ws.Range("BY9").PasteSpecial(-4122)
$ws.Range("BY10").Value = 147.6
$ws.Range("B2").Copy()
$ws.Range("BY10").PasteSpecial(-4122)
$ws.Range("BY11").Value = 201.3
$ws.Range("B2").Copy()
$ws.Range("BY11").PasteSpecial(-4122)
$ws.Range("BY12").Value = 163.8
$ws.Range("B2").Copy()
$ws.Range("BY12").PasteSpecial(-4122)
$ws.Range("BY13").Value = 175.2
$ws.Range("B2").Copy()
$ws.Range("BY13").PasteSpecial(-4122)
$ws.Range("BY14").Value = 141.1
$ws.Range("B2").Copy()
$ws.Range("BY14").PasteSpecial(-4122)
$ws.Range("BY15").Value = 196.2
$ws.Range("B2").Copy()
$ws.Range("BY15").PasteSpecial(-4122)
$ws.Range("BY16").Value = 124.3
$ws.Range("D2").Copy()
$ws.Range("BY16").PasteSpecial(-4122)
$ws.Range("BY17").Value = 207.1
$ws.Range("B2").Copy()
$ws.Range("BY17").PasteSpecial(-4122)
$ws.Range("BY18").Value = 170.9
$ws.Range("B2").Copy()
$ws.Range("BY18").PasteSpecial(-4122)
$ws.Range("BY19").Value = 123.4
$ws.Range("D2").Copy()
$ws.Range("BY19").PasteSpecial(-4122)
$ws.Range("BY20").Value = 149.9
$ws.Range("B2").Copy()
$ws.Range("BY20").PasteSpecial(-4122)
$ws.Range("BY21").Value = 215.3
$ws.Range("B2").Copy()
$ws.Range("BY21").PasteSpecial(-4122)
$ws.Range("BY22").Value = 125.8
$ws.Range("N2").Copy()
$ws.Range("BY22").PasteSpecial(-4122)
$ws.Range("BY23").Value = 133.1
$ws.Range("N2").Copy()
$ws.Range("BY23").PasteSpecial(-4122)
$ws.Range("BY24").Value = 128.9
$ws.Range("N2").Copy()
$ws.Range("BY24").PasteSpecial(-4122)
$ws.Range("BY25").Value = 124.9
$ws.Range("D2").Copy()
$ws.Range("BY25").PasteSpecial(-4122)
$ws.Range("BY26").Value = 164.2
$ws.Range("B2").Copy()
$ws.Range("BY26").PasteSpecial(-4122)
$ws.Range("BY27").Value = 186.1
$ws.Range("B2").Copy()
$ws.Range("BY27").PasteSpecial(-4122)
$ws.Range("BY28").Value = 163
$ws.Range("B2").Copy()
$ws.Range("BY28").PasteSpecial(-4122)
$ws.Range("BY29").Value = 237.6
$ws.Range("B2").Copy()
$ws.Range("BY29").PasteSpecial(-4122)
$ws.Range("BY30").Value = 175.4
$ws.Range("B2").Copy()
$ws.Range("BY30").PasteSpecial(-4122)
$ws.Range("BY31").Value = 175.2
$ws.Range("B2").Copy()
$ws.Range("BY31").PasteSpecial(-4122)
$ws.Range("BY32").Value = 133.2
$ws.Range("N2").Copy()
$ws.Range("BY32").PasteSpecial(-4122)
$ws.Range("BY33").Value = 144.4
$ws.Range("B2").Copy()
$ws.Range("BY33").PasteSpecial(-4122)
$ws.Range("BY34").Value = 156.2
$ws.Range("B2").Copy()
$ws.Range("BY34").PasteSpecial(-4122)
$ws.Range("BY35").Value = 113
$ws.Range("D2").Copy()
$ws.Range("BY35").PasteSpecial(-4122)
$ws.Range("BY36").Value = 181.1
$ws.Range("B2").Copy()
$ws.Range("BY36").PasteSpecial(-4122)
$ws.Range("BY37").Value = 168.9
$ws.Range("B2").Copy()
$ws.Range("BY37").PasteSpecial(-4122)
$ws.Range("BY38").Value = 152.7
$ws.Range("B2").Copy()
$ws.Range("BY38").PasteSpecial(-4122)
$ws.Range("BY39").Value = 185.1
$ws.Range("B2").Copy()
$ws.Range("BY39").PasteSpecial(-4122)
$ws.Range("BY40").Value = 232.8
$ws.Range("B2").Copy()
$ws.Range("BY40").PasteSpecial(-4122)
$ws.Range("BY41").Value = 136.6
$ws.Range("N2").Copy()
$ws.Range("BY41").PasteSpecial(-4122)
$ws.Range("BY42").Value = 151.6
$ws.Range("B2").Copy()
$ws.Range("BY42").PasteSpecial(-4122)
$ws.Range("BY43").Value = 140.1
$ws.Range("B2").Copy()
$ws.Range("BY43").PasteSpecial(-4122)
$ws.Range("BY44").Value = 139.6
$ws.Range("N2").Copy()
$ws.Range("BY44").PasteSpecial(-4122)
$ws.Range("BY45").Value = 180.3
$ws.Range("B2").Copy()
$ws.Range("BY45").PasteSpecial(-4122)
$ws.Range("BY46").Value = 149
$ws.Range("B2").Copy()
$ws.Range("BY46").PasteSpecial(-4122)
$ws.Range("BY47").Value = 164.2
$ws.Range("B2").Copy()
$ws.Range("BY47").PasteSpecial(-4122)
$ws.Range("BY48").Value = 181.7
$ws.Range("B2").Copy()
$ws.Range("BY48").PasteSpecial(-4122)
$ws.Range("BY49").Value = 163
$ws.Range("B2").Copy()
$ws.Range("BY49").PasteSpecial(-4122)
$ws.Range("BY50").Value = 111.8
$ws.Range("D2").Copy()
$ws.Range("BY50").PasteSpecial(-4122)
$ws.Range("BY51").Value = 139.2
$ws.Range("N2").Copy()
$ws.Range("BY51").PasteSpecial(-4122)
$ws.Range("BY52").Value = 122.9
$ws.Range("D2").Copy()
$ws.Range("BY52").PasteSpecial(-4122)
$ws.Range("BY53").Value = 159.9
$ws.Range("B2").Copy()
$ws.Range("BY53").PasteSpecial(-4122)
